$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "thej" / Developer task (previously row 4's data, now first row) ---
$ws.Range("A2").Value = '"5f0553090269c12d1c8ee6e7"'
$ws.Range("B2").Value = "thej"
$ws.Range("C2").Value = "Developer"
$ws.Range("D2").Value = "Pending"

# --- Row 3: new "Social Media" task, Done ---
$ws.Range("A3").Value = '"5f05accff6b3bb22dcc17178"'
$ws.Range("B3").Value = "thej"
$ws.Range("C3").Value = "Social Media"
$ws.Range("D3").Value = "Done"

# --- Row 4: new Efehi / Developer task, Pending ---
$ws.Range("A4").Value = '"5f05af8104dfeb226c3c42e0"'
$ws.Range("B4").Value = "Efehi"
$ws.Range("C4").Value = "Developer"
$ws.Range("D4").Value = "Pending"

# --- Column widths: all four columns become a uniform 100-characters wide ---
$ws.Columns.Item(1).ColumnWidth = 99.1
$ws.Columns.Item(2).ColumnWidth = 99.1
$ws.Columns.Item(3).ColumnWidth = 99.1
$ws.Columns.Item(4).ColumnWidth = 99.1

# Column D previously belonged to an outline group (outlineLevel=1, collapsed);
# clear the grouping so it no longer carries that old metadata.
$ws.Columns.Item(4).OutlineLevel = 0
